# Add a "Border Color" field to the Reports table, between "Font Name" and
# "Header Background Color". This shifts the existing H:N columns one slot
# to the right (into I:O) and fills in the new column H with the
# "Border Color" / "#333333" pair.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Columns")
$ws2 = $wb.Worksheets.Item("Reports")

# Shift the trailing columns (Header Background Color .. Print Pages Width)
# one column to the right, working from the rightmost column inward so we
# never overwrite a source cell before it has been read.
$srcCols = @("N", "M", "L", "K", "J", "I", "H")
$dstCols = @("O", "N", "M", "L", "K", "J", "I")

for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $src = $srcCols[$i]
    $dst = $dstCols[$i]
    $ws2.Range($dst + "1").Value = $ws2.Range($src + "1").Value2
    $ws2.Range($dst + "2").Value = $ws2.Range($src + "2").Value2
}

# New "Border Color" column
$ws2.Range("H1").Value = "Border Color"
$ws2.Range("H2").Value = "#333333"

# Grow the table to cover the new column, then make sure the last column's
# header name is picked up from the worksheet cell.
$lo = $ws2.ListObjects.Item("Table2")
$lo.Resize($ws2.Range("A1:O2"))
$ws2.Range("O1").Value = "Print Pages Width"

# Approximate column widths for the shifted/new columns (H:O). The new
# column H inherits the old column G's width, and every later column keeps
# its previous width shifted one slot to the right.
$widths = @{
    "H" = 11.4609375
    "I" = 23.15234375
    "J" = 17.23046875
    "K" = 16.23046875
    "L" = 16.4609375
    "M" = 17.765625
    "N" = 17
    "O" = 17.84375
}
foreach ($col in $widths.Keys) {
    $ws2.Range($col + "1").ColumnWidth = $widths[$col] - 0.8333333333333334
}

# Selection / active-sheet bookkeeping: the author finished on the
# "Columns" sheet at J4, then moved to "Reports" and ended on H3.
[void]$ws1.Range("J4").Select()
[void]$ws2.Activate()
[void]$ws2.Range("H3").Select()
